$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product title in A2 from "test product" to "test product5"
$ws.Range("A2").Value = "test product5"

# Update the active selection to match the saved view state (B11)
$ws.Range("B11").Select()
